# Fruta / hortaliza, semanal
# Insert 6 new weekly price rows for "Clemenuless" mandarina before the
# existing row 362, shifting the old rows 362-370 down to 368-376.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows at 362 (shifts existing rows 362-370 down to 368-376,
# inheriting formatting -- e.g. the date style -- from the former row 362).
$ws.Range("A362:A367").EntireRow.Insert(-4121)

# Common boilerplate values shared by every row in this block.
$mercadoId   = 2
$mercado     = "Comercializadora del Agro de Limarí"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102004
$categoria   = "Mandarina"
$origen      = "Provincia de Limarí"

function Set-Row {
    param($row, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidad, $precioKg, $kgUnidad)

    $ws.Range("A$row").Value = $mercadoId
    $ws.Range("B$row").Value = $mercado
    $ws.Range("C$row").Value = $region
    $ws.Range("D$row").Value = $fecha
    $ws.Range("E$row").Value = $codreg
    $ws.Range("F$row").Value = $tipo
    $ws.Range("G$row").Value = $productoId
    $ws.Range("H$row").Value = $producto
    $ws.Range("I$row").Value = $categoriaId
    $ws.Range("J$row").Value = $categoria
    $ws.Range("K$row").Value = $variedad
    $ws.Range("L$row").Value = $calidad
    $ws.Range("M$row").Value = $volumen
    $ws.Range("N$row").Value = $precioMin
    $ws.Range("O$row").Value = $precioMax
    $ws.Range("P$row").Value = $precioProm
    $ws.Range("Q$row").Value = $unidad
    $ws.Range("R$row").Value = $origen
    $ws.Range("S$row").Value = $precioKg
    $ws.Range("T$row").Value = $kgUnidad
}

$bandeja = "`$/bandeja 10 kilos"
$bins    = "`$/bins (450 kilos)"

Set-Row 362 44706 "Clemenuless" "Especial" 500 7500   8000   7750   $bandeja 775 10
Set-Row 363 44706 "Clemenuless" "Especial" 20  245000 250000 247500 $bins    550 450
Set-Row 364 44706 "Clemenuless" "Primera"  500 5500   6000   5750   $bandeja 575 10
Set-Row 365 44706 "Clemenuless" "Primera"  20  195000 200000 197500 $bins    439 450
Set-Row 366 44706 "Clemenuless" "Segunda"  400 3500   4000   3750   $bandeja 375 10
Set-Row 367 44706 "Clemenuless" "Segunda"  16  145000 150000 147500 $bins    328 450

Write-Output "Rows 362-367 inserted; sheet now has $($ws.UsedRange.Rows.Count) rows."
